$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the missing date/time entries for row 8 (security problems list entry)
# Date 2025-11-25 as Excel serial date, and time 20:16 as Excel time fraction
$ws.Range("B8").Value = 45986
$ws.Range("C8").Value = 0.84444444444444444

# Move the active selection to E5 (as left by the editor)
$ws.Range("E5").Select()
